# Update KDI simulation scripts
# - add biofuel simulation
# - commodity shock collection now includes source of energy

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Gasoline own-price elasticity figure (row 26) and add its source note
$ws.Range("F26").Value = -0.5755
$ws.Range("G26").Value = "average elasticity is -0.5755, max is -1.05, min is -0.16"
$ws.Range("I26").Value = "Dahl and Serner"

# Update the "Bernstein and Griffin" note in G14 to include confidence values
$ws.Range("G14").Value = "average residential elasticity is -0.2811 (0.0753), commercial elasticity is -0.96044 (0.627)"

# Update the view state: scroll so D4 is the top-left cell and select N12
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("N12").Select()
